# More CASP data formatting.
# Fill in the newly-recorded columns (L:V) for a batch of rows on Sheet1
# (ClosestDistance, AngleOfClosestDistance, Strongest_behavior, Sex, and the
# individual behavior-response minute columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("L132").Value = 30
$ws.Range("M132").Value = 180
$ws.Range("N132").Value = 2
$ws.Range("O132").Value = "M"
$ws.Range("Q132").Value = 30

$ws.Range("L135").Value = 20
$ws.Range("M135").Value = 265
$ws.Range("N135").Value = 2
$ws.Range("O135").Value = "M"

$ws.Range("L136").Value = 30
$ws.Range("M136").Value = 70
$ws.Range("N136").Value = 5
$ws.Range("O136").Value = "M"
$ws.Range("Q136").Value = 30
$ws.Range("T136").Value = 30

$ws.Range("L138").Value = "NA"
$ws.Range("M138").Value = "NA"
$ws.Range("N138").Value = 3
$ws.Range("O138").Value = "M"
$ws.Range("R138").Value = "NA"

$ws.Range("L139").Value = "NA"
$ws.Range("M139").Value = "NA"
$ws.Range("N139").Value = 5
$ws.Range("O139").Value = "M"
$ws.Range("R139").Value = "NA"
$ws.Range("T139").Value = "NA"

$ws.Range("L143").Value = 20
$ws.Range("M143").Value = 180
$ws.Range("N143").Value = 2
$ws.Range("O143").Value = "M"
$ws.Range("Q143").Value = 20

$ws.Range("L144").Value = 20
$ws.Range("M144").Value = 90
$ws.Range("N144").Value = 7
$ws.Range("O144").Value = "M"
$ws.Range("Q144").Value = 30
$ws.Range("V144").Value = 20

$ws.Range("L146").Value = 30
$ws.Range("M146").Value = 90
$ws.Range("N146").Value = 7
$ws.Range("O146").Value = "M"
$ws.Range("Q146").Value = 30
$ws.Range("R146").Value = 30
$ws.Range("V146").Value = 30

$ws.Range("L147").Value = 20
$ws.Range("M147").Value = 10
$ws.Range("N147").Value = 6
$ws.Range("O147").Value = "M"
$ws.Range("Q147").Value = 100
$ws.Range("U147").Value = 20

$ws.Range("L148").Value = 40
$ws.Range("M148").Value = 300
$ws.Range("N148").Value = 7
$ws.Range("O148").Value = "M"
$ws.Range("Q148").Value = 40
$ws.Range("V148").Value = 40

$ws.Range("L149").Value = 20
$ws.Range("M149").Value = 180
$ws.Range("N149").Value = 7
$ws.Range("O149").Value = "PAIR"
$ws.Range("Q149").Value = 20
$ws.Range("U149").Value = 20
$ws.Range("V149").Value = 20

$ws.Range("L150").Value = 10
$ws.Range("M150").Value = 70
$ws.Range("N150").Value = 2
$ws.Range("O150").Value = "M"
$ws.Range("Q150").Value = 10

$ws.Range("L152").Value = 5
$ws.Range("M152").Value = 245
$ws.Range("N152").Value = 5
$ws.Range("O152").Value = "M"
$ws.Range("Q152").Value = 30
$ws.Range("T152").Value = 5

$ws.Range("L153").Value = 10
$ws.Range("M153").Value = 320
$ws.Range("N153").Value = 7
$ws.Range("O153").Value = "PAIR"
$ws.Range("R153").Value = 60
$ws.Range("V153").Value = 60

$ws.Range("L154").Value = 5
$ws.Range("M154").Value = 95
$ws.Range("N154").Value = 7
$ws.Range("O154").Value = "M"
$ws.Range("Q154").Value = 50
$ws.Range("T154").Value = 5
$ws.Range("V154").Value = 5

$ws.Range("L155").Value = 30
$ws.Range("M155").Value = 270
$ws.Range("N155").Value = 6
$ws.Range("O155").Value = "M"
$ws.Range("Q155").Value = 30
$ws.Range("U155").Value = 30

$ws.Range("L157").Value = 20
$ws.Range("M157").Value = 0
$ws.Range("N157").Value = 2
$ws.Range("O157").Value = "M"

$ws.Range("L158").Value = 60
$ws.Range("M158").Value = 275
$ws.Range("N158").Value = 7
$ws.Range("O158").Value = "M"
$ws.Range("Q158").Value = 60
$ws.Range("V158").Value = 60

$ws.Range("L159").Value = 40
$ws.Range("M159").Value = 190
$ws.Range("N159").Value = 7
$ws.Range("O159").Value = "M"
$ws.Range("V159").Value = 40

$ws.Range("L160").Value = 10
$ws.Range("M160").Value = 180
$ws.Range("N160").Value = 7
$ws.Range("O160").Value = "M"
$ws.Range("Q160").Value = 10
$ws.Range("S160").Value = 10
$ws.Range("V160").Value = 10

$ws.Range("L161").Value = 30
$ws.Range("M161").Value = 190
$ws.Range("N161").Value = 2
$ws.Range("O161").Value = "UNK"
$ws.Range("Q161").Value = 30

$ws.Range("L163").Value = 20
$ws.Range("M163").Value = 265
$ws.Range("N163").Value = 3
$ws.Range("O163").Value = "M"
$ws.Range("Q163").Value = 20
$ws.Range("R163").Value = 40

$ws.Range("L165").Value = 40
$ws.Range("M165").Value = 90
$ws.Range("N165").Value = 4
$ws.Range("O165").Value = "M"
$ws.Range("Q165").Value = 40
$ws.Range("S165").Value = 40

$ws.Range("L175").Value = 50
$ws.Range("M175").Value = 90
$ws.Range("N175").Value = 4
$ws.Range("O175").Value = "M"
$ws.Range("S175").Value = 50

$ws.Range("L181").Value = 40
$ws.Range("M181").Value = 180
$ws.Range("N181").Value = 2
$ws.Range("O181").Value = "UNK"
$ws.Range("Q181").Value = 40

$ws.Range("L182").Value = 100
$ws.Range("M182").Value = 180
$ws.Range("N182").Value = 7
$ws.Range("O182").Value = "M"
$ws.Range("V182").Value = 100

# Update the frozen-pane scroll position and the active selection in the
# bottom-right pane to reflect where the user was working after the edits.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 248
$win.ScrollColumn = 11
$ws.Range("L183").Select()
